$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.336.31'
$ws.Range("E2").Value = '  -4.77%  '
$ws.Range("D3").Value = '3.093.98'
$ws.Range("E3").Value = '  -4.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.59'
$ws.Range("E5").Value = '  -4.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.41'
$ws.Range("E6").Value = '  -11.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.089.92'
$ws.Range("E8").Value = '  -4.27%  '
$ws.Range("E9").Value = '  -3.83%  '
$ws.Range("E10").Value = '  -5.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.34'
$ws.Range("E11").Value = '  -10.13%  '
$ws.Range("E12").Value = '  -4.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.45'
$ws.Range("E13").Value = '  -6.63%  '
$ws.Range("E14").Value = '  -7.81%  '
$ws.Range("D15").Value = '3.593.43'
$ws.Range("E15").Value = '  -4.40%  '
$ws.Range("D16").Value = '63.390.29'
$ws.Range("E16").Value = '  -4.64%  '
$ws.Range("E17").Value = '  -3.39%  '
$ws.Range("D18").Value = '3.095.09'
$ws.Range("E18").Value = '  -4.52%  '
$ws.Range("E19").Value = '  -5.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '491.24'
$ws.Range("E20").Value = '  -12.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.67'
$ws.Range("E21").Value = '  -5.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.716'
$ws.Range("E22").Value = '  -3.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.25'
$ws.Range("E23").Value = '  -7.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.01'
$ws.Range("E24").Value = '  -3.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.41'
$ws.Range("E25").Value = '  -8.83%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.47'
$ws.Range("E27").Value = '  -10.26%  '
$ws.Range("E28").Value = '  -6.53%  '
$ws.Range("E29").Value = '  -12.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.61'
$ws.Range("E31").Value = '  -4.55%  '
$ws.Range("E32").Value = '  -3.75%  '
$ws.Range("E33").Value = '  -9.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '58.88'
$ws.Range("E34").Value = '  +5.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '520.95'
$ws.Range("E35").Value = '  -8.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.01'
$ws.Range("E36").Value = '  -6.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.16'
$ws.Range("E37").Value = '  -10.93%  '
$ws.Range("E38").Value = '  -11.62%  '
$ws.Range("D39").Value = '3.143.02'
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("E40").Value = '  -7.65%  '
$ws.Range("E41").Value = '  -5.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.18'
$ws.Range("E42").Value = '  -5.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.68'
$ws.Range("E43").Value = '  -12.46%  '
$ws.Range("E44").Value = '  -6.14%  '
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("E46").Value = '  -11.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.16'
$ws.Range("E47").Value = '  -6.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '121.64'
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("E49").Value = '  -4.24%  '
$ws.Range("E50").Value = '  -10.47%  '
$ws.Range("E51").Value = '  -9.82%  '
